$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Activate the "Бэклог задач" sheet/window so it becomes the active tab
# (tabSelected on sheet2, activeTab on the workbook, and tabSelected removed
# from sheet1 all follow from this).
$ws.Activate()

# Append four new backlog rows (42-45), column B = text, column C = timestamp.
$ws.Cells.Item(42, 2).Value = "Добавить обновление rating.updated, при обновлении позиции."
$ws.Cells.Item(42, 3).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(42, 3).Value = (Get-Date -Year 2015 -Month 1 -Day 19 -Hour 17 -Minute 14 -Second 0)

$ws.Cells.Item(43, 2).Value = "Сделать логирование по id класса! "
$ws.Cells.Item(43, 3).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(43, 3).Value = (Get-Date -Year 2015 -Month 1 -Day 19 -Hour 17 -Minute 17 -Second 0)

$ws.Cells.Item(44, 2).Value = "Сделать логирование по id сообщения!"
$ws.Cells.Item(44, 3).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(44, 3).Value = (Get-Date -Year 2015 -Month 1 -Day 19 -Hour 17 -Minute 17 -Second 0)

$ws.Cells.Item(45, 2).Value = "Организовать событийнные механизмы. Ато как то евент вызыаются прямыми обращениями в методы."
$ws.Cells.Item(45, 3).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(45, 3).Value = (Get-Date -Year 2015 -Month 1 -Day 19 -Hour 17 -Minute 24 -Second 0)

# Move the selection to the new last row, matching the recorded cursor position.
$ws.Range("C46").Select()
